# Add new survey respondent rows to sheet1 (Когтевран), sheet2 (Гриффиндор)
# and sheet3 (Слизерин), matching the "checking process for name persons"
# update: new rows of student responses appended under the existing header
# + first data row on each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Когтевран -- add rows 3..6
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Когтевран")

$ws1.Range("A3").Value = 2
$ws1.Range("B3").Value = "Анна"
$ws1.Range("C3").Value = 3
$ws1.Range("D3").Value = 3
$ws1.Range("E3").Value = 3
$ws1.Range("F3").Value = 0
$ws1.Range("G3").Value = "https://vk.com/id739914548"
$ws1.Range("H3").Value = 44758.78654422454

$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = "Олеся Грейнджер"
$ws1.Range("C4").Value = 3
$ws1.Range("D4").Value = 2
$ws1.Range("E4").Value = 2
$ws1.Range("F4").Value = 2
$ws1.Range("G4").Value = "https://vk.com/id445175392"
$ws1.Range("H4").Value = 44758.81101667824

$ws1.Range("A5").Value = 4
$ws1.Range("B5").Value = "Алекса Черни"
$ws1.Range("C5").Value = 4
$ws1.Range("D5").Value = 2
$ws1.Range("E5").Value = 0
$ws1.Range("F5").Value = 3
$ws1.Range("G5").Value = "https://vk.com/id576527766"
$ws1.Range("H5").Value = 44759.35237075231

$ws1.Range("A6").Value = 5
$ws1.Range("B6").Value = "Иренчик"
$ws1.Range("C6").Value = 6
$ws1.Range("D6").Value = 1
$ws1.Range("E6").Value = 0
$ws1.Range("F6").Value = 2
$ws1.Range("G6").Value = "https://vk.com/id332431318"
$ws1.Range("H6").Value = 44759.98855631945

$ws1.Range("H2").Copy()
$ws1.Range("H3:H6").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Sheet 2: Гриффиндор -- add rows 3..4
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Гриффиндор")

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = "Делисия Ви-Марет"
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 0
$ws2.Range("E3").Value = 5
$ws2.Range("F3").Value = 3
$ws2.Range("G3").Value = "https://vk.com/id166767722"
$ws2.Range("H3").Value = 44758.81439315972

$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = "Мак Так"
$ws2.Range("C4").Value = 3
$ws2.Range("D4").Value = 0
$ws2.Range("E4").Value = 6
$ws2.Range("F4").Value = 0
$ws2.Range("G4").Value = "https://vk.com/id435298354"
$ws2.Range("H4").Value = 44759.710700891206

$ws2.Range("H2").Copy()
$ws2.Range("H3:H4").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Sheet 3: Слизерин -- add rows 2..4
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Слизерин")

$ws3.Range("A2").Value = 1
$ws3.Range("B2").Value = "Драко Малфой"
$ws3.Range("C2").Value = 3
$ws3.Range("D2").Value = 0
$ws3.Range("E2").Value = 3
$ws3.Range("F2").Value = 3
$ws3.Range("G2").Value = "https://vk.com/id710936448"
$ws3.Range("H2").Value = 44758.80729621528

$ws3.Range("A3").Value = 2
$ws3.Range("B3").Value = "Alhajia Unuk"
$ws3.Range("C3").Value = 1
$ws3.Range("D3").Value = 1
$ws3.Range("E3").Value = 3
$ws3.Range("F3").Value = 4
$ws3.Range("G3").Value = "https://vk.com/id733827532"
$ws3.Range("H3").Value = 44758.935105150464

$ws3.Range("A4").Value = 3
$ws3.Range("B4").Value = "Аделия Реддл"
$ws3.Range("C4").Value = 2
$ws3.Range("D4").Value = 2
$ws3.Range("E4").Value = 2
$ws3.Range("F4").Value = 3
$ws3.Range("G4").Value = "https://vk.com/id391472095"
$ws3.Range("H4").Value = 44759.73308201389

# Sheet 3 had no pre-existing data row w/ the date style, so copy the
# date-formatted style from sheet1's H2 (same workbook style used
# throughout for the registration-date column).
$ws1.Range("H2").Copy()
$ws3.Range("H2:H4").PasteSpecial(-4122)
